$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily records appended after row 238 (through 02/05/2021), as per
# commit "aggiornamento fino a 02/05".
$rows = @(
    @{ Row = 239; A = 44313; B = 0; C = 2; D = 81.59934720522236 },
    @{ Row = 240; A = 44314; B = 0; C = 1; D = 40.79967360261118 },
    @{ Row = 241; A = 44315; B = 0; C = 0; D = 0 },
    @{ Row = 242; A = 44316; B = 1; C = 1; D = 40.79967360261118 },
    @{ Row = 243; A = 44317; B = 0; C = 1; D = 40.79967360261118 },
    @{ Row = 244; A = 44318; B = 0; C = 1; D = 40.79967360261118 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Copy the formatting of the date cell above (style index "2": bordered,
    # bold, centered, datetime number format) onto the new date cell so no
    # new style entries are minted.
    $ws.Range("A" + ($rowNum - 1)).Copy()
    $ws.Range("A" + $rowNum).PasteSpecial(-4122)

    $ws.Range("A" + $rowNum).Value = $r.A
    $ws.Range("B" + $rowNum).Value = $r.B
    $ws.Range("C" + $rowNum).Value = $r.C
    $ws.Range("D" + $rowNum).Value = $r.D
}
